# Applies updated odds values to specific cells in the active worksheet,
# matching the targeted data refresh for Jogos_da_Semana_FlashScore_2025-06-07.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 7).Value = 2.9
$ws.Cells.Item(2, 9).Value = 2.9
$ws.Cells.Item(2, 10).Value = 1.18
$ws.Cells.Item(2, 11).Value = 4.5
$ws.Cells.Item(2, 21).Value = 12
$ws.Cells.Item(2, 32).Value = 13

# Row 5
$ws.Cells.Item(5, 7).Value = 3.1
$ws.Cells.Item(5, 9).Value = 2.37
$ws.Cells.Item(5, 11).Value = 6.4
$ws.Cells.Item(5, 13).Value = 2.72
$ws.Cells.Item(5, 20).Value = 8.25
$ws.Cells.Item(5, 21).Value = 16.5
$ws.Cells.Item(5, 26).Value = 6.4
$ws.Cells.Item(5, 27).Value = 6.2
$ws.Cells.Item(5, 30).Value = 6.6
$ws.Cells.Item(5, 31).Value = 11.25
$ws.Cells.Item(5, 33).Value = 27
$ws.Cells.Item(5, 36).Value = 900

# Row 6
$ws.Cells.Item(6, 7).Value = 2.55
$ws.Cells.Item(6, 8).Value = 2.88
$ws.Cells.Item(6, 9).Value = 3.1
$ws.Cells.Item(6, 10).Value = 1.17
$ws.Cells.Item(6, 11).Value = 5
$ws.Cells.Item(6, 18).Value = 2.25
$ws.Cells.Item(6, 19).Value = 1.57
$ws.Cells.Item(6, 22).Value = 12
$ws.Cells.Item(6, 26).Value = 5

# Row 7
$ws.Cells.Item(7, 7).Value = 2.92
$ws.Cells.Item(7, 9).Value = 2.35
$ws.Cells.Item(7, 14).Value = 2.32
$ws.Cells.Item(7, 18).Value = 2.02
$ws.Cells.Item(7, 20).Value = 7
$ws.Cells.Item(7, 21).Value = 13
$ws.Cells.Item(7, 22).Value = 11.5
$ws.Cells.Item(7, 23).Value = 35
$ws.Cells.Item(7, 31).Value = 10
$ws.Cells.Item(7, 33).Value = 24
$ws.Cells.Item(7, 34).Value = 24

# Row 10
$ws.Cells.Item(10, 7).Value = 2.1
$ws.Cells.Item(10, 9).Value = 3.5
$ws.Cells.Item(10, 12).Value = 1.5
$ws.Cells.Item(10, 13).Value = 2.5
$ws.Cells.Item(10, 16).Value = 1.53
$ws.Cells.Item(10, 17).Value = 2.38
$ws.Cells.Item(10, 18).Value = 2.2
$ws.Cells.Item(10, 19).Value = 1.62
$ws.Cells.Item(10, 26).Value = 7
$ws.Cells.Item(10, 34).Value = 34
$ws.Cells.Item(10, 35).Value = 41

# Row 11
$ws.Cells.Item(11, 8).Value = 3.95
$ws.Cells.Item(11, 9).Value = 4.4
$ws.Cells.Item(11, 13).Value = 3.7
$ws.Cells.Item(11, 18).Value = 1.62
$ws.Cells.Item(11, 34).Value = 37

# Row 15
$ws.Cells.Item(15, 7).Value = 2.32
$ws.Cells.Item(15, 12).Value = 1.44
$ws.Cells.Item(15, 13).Value = 2.4
$ws.Cells.Item(15, 14).Value = 2.27
$ws.Cells.Item(15, 15).Value = 1.5
$ws.Cells.Item(15, 20).Value = 6.3
$ws.Cells.Item(15, 24).Value = 23
$ws.Cells.Item(15, 26).Value = 7.2
$ws.Cells.Item(15, 28).Value = 18
$ws.Cells.Item(15, 30).Value = 7.3
$ws.Cells.Item(15, 32).Value = 11.5
$ws.Cells.Item(15, 35).Value = 50

# Row 17
$ws.Cells.Item(17, 10).Value = 1.06
$ws.Cells.Item(17, 11).Value = 10
$ws.Cells.Item(17, 12).Value = 1.33
$ws.Cells.Item(17, 13).Value = 3.25
$ws.Cells.Item(17, 14).Value = 2.08
$ws.Cells.Item(17, 15).Value = 1.73

# Row 19
$ws.Cells.Item(19, 9).Value = 2.3
$ws.Cells.Item(19, 22).Value = 11
$ws.Cells.Item(19, 30).Value = 10
$ws.Cells.Item(19, 34).Value = 17

# Row 22
$ws.Cells.Item(22, 16).Value = 1.47
$ws.Cells.Item(22, 17).Value = 2.52
$ws.Cells.Item(22, 19).Value = 1.8
$ws.Cells.Item(22, 32).Value = 14

# Row 26
$ws.Cells.Item(26, 10).Value = 1.07
$ws.Cells.Item(26, 11).Value = 9
$ws.Cells.Item(26, 14).Value = 2.15
$ws.Cells.Item(26, 15).Value = 1.67

# Row 29
$ws.Cells.Item(29, 7).Value = 2.32
$ws.Cells.Item(29, 9).Value = 2.9
$ws.Cells.Item(29, 12).Value = 1.33
$ws.Cells.Item(29, 13).Value = 2.8
$ws.Cells.Item(29, 14).Value = 1.98
$ws.Cells.Item(29, 15).Value = 1.65
$ws.Cells.Item(29, 16).Value = 1.44
$ws.Cells.Item(29, 17).Value = 2.42
$ws.Cells.Item(29, 19).Value = 1.85
$ws.Cells.Item(29, 20).Value = 7.4
$ws.Cells.Item(29, 21).Value = 11
$ws.Cells.Item(29, 25).Value = 32
$ws.Cells.Item(29, 26).Value = 8.75
$ws.Cells.Item(29, 28).Value = 14.5
$ws.Cells.Item(29, 31).Value = 14.5
$ws.Cells.Item(29, 32).Value = 10.5
$ws.Cells.Item(29, 34).Value = 26

# Row 30
$ws.Cells.Item(30, 7).Value = 2.9
$ws.Cells.Item(30, 9).Value = 2.3
$ws.Cells.Item(30, 21).Value = 17
$ws.Cells.Item(30, 26).Value = 13
$ws.Cells.Item(30, 36).Value = 151

# Row 33
$ws.Cells.Item(33, 8).Value = 3.1
$ws.Cells.Item(33, 20).Value = 5.6
$ws.Cells.Item(33, 28).Value = 19.5
$ws.Cells.Item(33, 31).Value = 18
$ws.Cells.Item(33, 32).Value = 14
$ws.Cells.Item(33, 35).Value = 65

# Row 34
$ws.Cells.Item(34, 7).Value = 1.39
$ws.Cells.Item(34, 8).Value = 4.5
$ws.Cells.Item(34, 9).Value = 6.3
$ws.Cells.Item(34, 10).Value = 1.03
$ws.Cells.Item(34, 11).Value = 9
$ws.Cells.Item(34, 12).Value = 1.2
$ws.Cells.Item(34, 13).Value = 4.1
$ws.Cells.Item(34, 14).Value = 1.6
$ws.Cells.Item(34, 15).Value = 2.2
$ws.Cells.Item(34, 16).Value = 1.31
$ws.Cells.Item(34, 17).Value = 3.15
$ws.Cells.Item(34, 18).Value = 1.85
$ws.Cells.Item(34, 19).Value = 1.87
$ws.Cells.Item(34, 20).Value = 7.8
$ws.Cells.Item(34, 21).Value = 7.1
$ws.Cells.Item(34, 22).Value = 8.25
$ws.Cells.Item(34, 23).Value = 9.25
$ws.Cells.Item(34, 24).Value = 11
$ws.Cells.Item(34, 25).Value = 25
$ws.Cells.Item(34, 26).Value = 9
$ws.Cells.Item(34, 27).Value = 9.25
$ws.Cells.Item(34, 28).Value = 18.5
$ws.Cells.Item(34, 29).Value = 80
$ws.Cells.Item(34, 30).Value = 19
$ws.Cells.Item(34, 32).Value = 21
$ws.Cells.Item(34, 34).Value = 65
$ws.Cells.Item(34, 35).Value = 60
$ws.Cells.Item(34, 36).Value = 600

# Row 35
$ws.Cells.Item(35, 8).Value = 3.6
$ws.Cells.Item(35, 12).Value = 1.28
$ws.Cells.Item(35, 13).Value = 3.35
$ws.Cells.Item(35, 17).Value = 2.75
$ws.Cells.Item(35, 19).Value = 1.88
$ws.Cells.Item(35, 27).Value = 7.2
$ws.Cells.Item(35, 28).Value = 16
$ws.Cells.Item(35, 29).Value = 75
$ws.Cells.Item(35, 30).Value = 14
$ws.Cells.Item(35, 31).Value = 30

# Row 38
$ws.Cells.Item(38, 23).Value = 22
$ws.Cells.Item(38, 27).Value = 6.6
$ws.Cells.Item(38, 28).Value = 14.5
$ws.Cells.Item(38, 30).Value = 9
$ws.Cells.Item(38, 31).Value = 14
$ws.Cells.Item(38, 32).Value = 10.5

# Row 42
$ws.Cells.Item(42, 8).Value = 3.45
$ws.Cells.Item(42, 12).Value = 1.29
$ws.Cells.Item(42, 13).Value = 3
$ws.Cells.Item(42, 14).Value = 1.85
$ws.Cells.Item(42, 15).Value = 1.75
$ws.Cells.Item(42, 16).Value = 1.4
$ws.Cells.Item(42, 17).Value = 2.55
$ws.Cells.Item(42, 20).Value = 7.5
$ws.Cells.Item(42, 28).Value = 15
$ws.Cells.Item(42, 31).Value = 16
$ws.Cells.Item(42, 34).Value = 28
$ws.Cells.Item(42, 35).Value = 37
